$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.183.38'
$ws.Range('E2').Value = '  -3.81%  '
$ws.Range('D3').Value = '3.396.68'
$ws.Range('E3').Value = '  -4.45%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'577.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.15%  '
$ws.Range('D6').Value = "'131.01"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.53%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.394.42'
$ws.Range('E9').Value = '  -7.65%  '
$ws.Range('E10').Value = '  -10.23%  '
$ws.Range('E11').Value = '  -10.75%  '
$ws.Range('D12').Value = "'0.371"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.28%  '
$ws.Range('D13').Value = '3.977.86'
$ws.Range('E13').Value = '  -4.38%  '
$ws.Range('D14').Value = "'0.0000175"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -10.93%  '
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '64.257.68'
$ws.Range('E16').Value = '  -3.57%  '
$ws.Range('D17').Value = "'25.77"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -11.32%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.351.38'
$ws.Range('E18').Value = '  -5.65%  '
$ws.Range('D19').Value = "'9.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -13.49%  '
$ws.Range('D20').Value = "'5.63"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.12%  '
$ws.Range('D21').Value = "'13.49"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.46%  '
$ws.Range('D22').Value = "'377.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -11.65%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  -10.19%  '
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = "'71.21"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.28%  '
$ws.Range('D27').Value = '3.532.03'
$ws.Range('D28').Value = "'0.0000102"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -12.26%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = "'7.05"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E31').Value = '  -12.83%  '
$ws.Range('D32').Value = "'7.92"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -12.47%  '
$ws.Range('D33').Value = '3.411.78'
$ws.Range('E33').Value = '  -4.23%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -6.87%  '
$ws.Range('E36').Value = '  -11.02%  '
$ws.Range('D37').Value = "'170.79"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.58%  '
$ws.Range('E38').Value = '  -14.41%  '
$ws.Range('D39').Value = "'6.61"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.74%  '
$ws.Range('D40').Value = "'1.43"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.68%  '
$ws.Range('D41').Value = "'4.56"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -14.23%  '
$ws.Range('E42').Value = '  -9.29%  '
$ws.Range('E43').Value = '  -8.50%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').Value = "'41.58"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.57%  '
$ws.Range('D46').Value = "'4.24"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -15.97%  '
$ws.Range('E47').Value = '  -11.61%  '
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('D49').Value = "'22.00"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.79%  '
$ws.Range('D51').Value = '2.177.81'
$ws.Range('E51').Value = '  -7.04%  '
